$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 20:06"

# Apply updated statistics / reordered country rows as per the refreshed data pull
# Row 4
$ws.Cells.Item(4, 2).Value = 2220485
$ws.Cells.Item(4, 3).Value = 12085
$ws.Cells.Item(4, 4).Value = 904182
$ws.Cells.Item(4, 5).Value = 1196839
$ws.Cells.Item(4, 7).Value = 332
$ws.Cells.Item(4, 8).Value = 119464

# Row 7
$ws.Cells.Item(7, 2).Value = 363209
$ws.Cells.Item(7, 3).Value = 9048
$ws.Cells.Item(7, 4).Value = 192594
$ws.Cells.Item(7, 5).Value = 158550
$ws.Cells.Item(7, 7).Value = 144
$ws.Cells.Item(7, 8).Value = 12065

# Row 14
$ws.Cells.Item(14, 5).Value = 24602
$ws.Cells.Item(14, 7).Value = 232
$ws.Cells.Item(14, 8).Value = 3615

# Row 15
$ws.Cells.Item(15, 2).Value = 182727
$ws.Cells.Item(15, 3).Value = 1429
$ws.Cells.Item(15, 4).Value = 154640
$ws.Cells.Item(15, 5).Value = 23226
$ws.Cells.Item(15, 7).Value = 19
$ws.Cells.Item(15, 8).Value = 4861

# Row 32
$ws.Cells.Item(32, 2).Value = 43364
$ws.Cells.Item(32, 3).Value = 382
$ws.Cells.Item(32, 4).Value = 29537
$ws.Cells.Item(32, 5).Value = 13532
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 8).Value = 295

# Row 44
$ws.Cells.Item(44, 2).Value = 25341
$ws.Cells.Item(44, 3).Value = 7
$ws.Cells.Item(44, 5).Value = 933
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 1710

# Row 50
$ws.Cells.Item(50, 2).Value = 19783
$ws.Cells.Item(50, 3).Value = 288
$ws.Cells.Item(50, 5).Value = 4021

# Row 56
$ws.Cells.Item(56, 5).Value = 5525
$ws.Cells.Item(56, 7).Value = 9
$ws.Cells.Item(56, 8).Value = 97

# Row 63
$ws.Cells.Item(63, 1).Value = "Azerbaiyan"
$ws.Cells.Item(63, 2).Value = 10991
$ws.Cells.Item(63, 3).Value = 329
$ws.Cells.Item(63, 4).Value = 6075
$ws.Cells.Item(63, 5).Value = 4783
$ws.Cells.Item(63, 7).Value = 7
$ws.Cells.Item(63, 8).Value = 133

# Row 64
$ws.Cells.Item(64, 1).Value = "Guatemala"
$ws.Cells.Item(64, 2).Value = 10706
$ws.Cells.Item(64, 3).Value = 434
$ws.Cells.Item(64, 4).Value = 2096
$ws.Cells.Item(64, 5).Value = 8192
$ws.Cells.Item(64, 7).Value = 19
$ws.Cells.Item(64, 8).Value = 418

# Row 68
$ws.Cells.Item(68, 2).Value = 8997
$ws.Cells.Item(68, 3).Value = 66
$ws.Cells.Item(68, 4).Value = 7993
$ws.Cells.Item(68, 5).Value = 791
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 213

# Row 101
$ws.Cells.Item(101, 2).Value = 2120
$ws.Cells.Item(101, 3).Value = 26
$ws.Cells.Item(101, 4).Value = 1677
$ws.Cells.Item(101, 5).Value = 435

# Row 121
$ws.Cells.Item(121, 2).Value = 1308
$ws.Cells.Item(121, 3).Value = 5
$ws.Cells.Item(121, 4).Value = 711
$ws.Cells.Item(121, 5).Value = 584

# Row 122
$ws.Cells.Item(122, 2).Value = 1249
$ws.Cells.Item(122, 3).Value = 24
$ws.Cells.Item(122, 4).Value = 707
$ws.Cells.Item(122, 5).Value = 491

# Row 129
$ws.Cells.Item(129, 1).Value = "Yemen"
$ws.Cells.Item(129, 2).Value = 902
$ws.Cells.Item(129, 3).Value = 17
$ws.Cells.Item(129, 4).Value = 271
$ws.Cells.Item(129, 5).Value = 387
$ws.Cells.Item(129, 7).Value = 30
$ws.Cells.Item(129, 8).Value = 244

# Row 130
$ws.Cells.Item(130, 1).Value = "Burkina Faso"
$ws.Cells.Item(130, 2).Value = 895
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 809
$ws.Cells.Item(130, 5).Value = 33
$ws.Cells.Item(130, 8).Value = 53

# Row 131
$ws.Cells.Item(131, 1).Value = "Georgia"
$ws.Cells.Item(131, 2).Value = 888
$ws.Cells.Item(131, 3).Value = 9
$ws.Cells.Item(131, 4).Value = 731
$ws.Cells.Item(131, 5).Value = 143
$ws.Cells.Item(131, 8).Value = 14

# Row 133
$ws.Cells.Item(133, 1).Value = "Republica del Chad"
$ws.Cells.Item(133, 3).Value = 1
$ws.Cells.Item(133, 4).Value = 721
$ws.Cells.Item(133, 5).Value = 59
$ws.Cells.Item(133, 8).Value = 74

# Row 134
$ws.Cells.Item(134, 1).Value = "Principado de Andorra"
$ws.Cells.Item(134, 2).Value = 854
$ws.Cells.Item(134, 4).Value = 791
$ws.Cells.Item(134, 5).Value = 11
$ws.Cells.Item(134, 8).Value = 52

# Row 147
$ws.Cells.Item(147, 1).Value = "Suazilandia"
$ws.Cells.Item(147, 2).Value = 563
$ws.Cells.Item(147, 3).Value = 43
$ws.Cells.Item(147, 4).Value = 262
$ws.Cells.Item(147, 5).Value = 297
$ws.Cells.Item(147, 8).Value = 4

# Row 148
$ws.Cells.Item(148, 1).Value = "Togo"
$ws.Cells.Item(148, 2).Value = 537
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 344
$ws.Cells.Item(148, 5).Value = 180
$ws.Cells.Item(148, 8).Value = 13

# Row 149
$ws.Cells.Item(149, 1).Value = "Estado de Palestina"
$ws.Cells.Item(149, 2).Value = 532
$ws.Cells.Item(149, 3).Value = 18
$ws.Cells.Item(149, 4).Value = 415
$ws.Cells.Item(149, 5).Value = 114
$ws.Cells.Item(149, 8).Value = 3

# Row 152
$ws.Cells.Item(152, 2).Value = 497
$ws.Cells.Item(152, 3).Value = 2
$ws.Cells.Item(152, 5).Value = 36

# Row 155
$ws.Cells.Item(155, 2).Value = 401
$ws.Cells.Item(155, 3).Value = 10
$ws.Cells.Item(155, 4).Value = 63
$ws.Cells.Item(155, 5).Value = 334

# Row 174
$ws.Cells.Item(174, 1).Value = "Eritrea"
$ws.Cells.Item(174, 2).Value = 131
$ws.Cells.Item(174, 3).Value = 10
$ws.Cells.Item(174, 4).Value = 39
$ws.Cells.Item(174, 5).Value = 92

# Row 175
$ws.Cells.Item(175, 1).Value = "Camboya"
$ws.Cells.Item(175, 2).Value = 128
$ws.Cells.Item(175, 4).Value = 126
$ws.Cells.Item(175, 5).Value = 2
$ws.Cells.Item(175, 8).Value = 0

# Row 176
$ws.Cells.Item(176, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(176, 2).Value = 123
$ws.Cells.Item(176, 4).Value = 109
$ws.Cells.Item(176, 5).Value = 6
$ws.Cells.Item(176, 8).Value = 8

# Row 183
$ws.Cells.Item(183, 1).Value = "Botsuana"
$ws.Cells.Item(183, 2).Value = 79
$ws.Cells.Item(183, 3).Value = 19
$ws.Cells.Item(183, 4).Value = 24
$ws.Cells.Item(183, 5).Value = 54
$ws.Cells.Item(183, 8).Value = 1

# Row 184
$ws.Cells.Item(184, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(184, 2).Value = 77
$ws.Cells.Item(184, 4).Value = 62
$ws.Cells.Item(184, 5).Value = 0
$ws.Cells.Item(184, 8).Value = 15

# Row 185
$ws.Cells.Item(185, 1).Value = "Butan"
$ws.Cells.Item(185, 2).Value = 67
$ws.Cells.Item(185, 5).Value = 43
$ws.Cells.Item(185, 8).Value = 0

# Row 190
$ws.Cells.Item(190, 2).Value = 36
$ws.Cells.Item(190, 3).Value = 2
$ws.Cells.Item(190, 4).Value = 19
$ws.Cells.Item(190, 5).Value = 17

# Row 206
$ws.Cells.Item(206, 1).Value = "Groenlandia"

# Row 207
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"

# Row 210
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

# Row 211
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1
